$d = $word.ActiveDocument

# The "Solutions" table is the first (only) table in the document.
$t = $d.Tables.Item(1)

# Append a new row for Problem 7, mirroring the formatting of the
# preceding rows (Word inherits pStyle/jc from the last row automatically).
$newRow = $t.Rows.Add()

$newRow.Cells.Item(1).Range.Text = "7"
$newRow.Cells.Item(2).Range.Text = "-"
$newRow.Cells.Item(3).Range.Text = "0.377"
